$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 310; this shifts the existing rows 310-336 down to 311-337
$ws.Rows.Item(310).Insert()

# Populate the newly inserted row 310 with its data
$ws.Range("A310").Value = 8
$ws.Range("B310").Value = "Terminal La Palmera de La Serena"
$ws.Range("C310").Value = "Coquimbo"
$ws.Range("D310").Value = 45106
$ws.Range("D310").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E310").Value = 4
$ws.Range("F310").Value = 100112037
$ws.Range("G310").Value = "Cebollín"
$ws.Range("H310").Value = "Sin especificar"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 900
$ws.Range("K310").Value = 1000
$ws.Range("L310").Value = 1200
$ws.Range("M310").Value = 1100
$ws.Range("N310").Value = '$/paquete 6 unidades'
$ws.Range("O310").Value = "Provincia del Elquí"
$ws.Range("P310").Value = 183
$ws.Range("Q310").Value = 6
$ws.Range("R310").Value = "Hortaliza"
